# Apply the 2022-11-29 11:12:18 re-crawl update to the Coop "bread" export.
#
# 1) A handful of product rows were re-ordered within the sheet (the
#    underlying scraper emitted them in a different sequence on the second
#    run). We reproduce that by swapping/rotating the full row content
#    (columns A:N) between the affected row numbers, leaving every other
#    row untouched.
# 2) Every data row's timestamp (column O) is bumped from
#    "2022-11-29 10:56:57" to "2022-11-29 11:12:18".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Get-RowData {
    param($ws, $row)
    $data = @{}
    for ($c = 1; $c -le 14; $c++) {
        $data[$c] = $ws.Cells.Item($row, $c).Value2
    }
    return $data
}

function Set-RowData {
    param($ws, $row, $data)
    for ($c = 1; $c -le 14; $c++) {
        $val = $data[$c]
        $cell = $ws.Cells.Item($row, $c)
        if ($c -eq 4 -or $c -eq 5) {
            # D = ratingAmount, E = ratingValue -> numeric columns
            if ($null -eq $val) {
                $cell.ClearContents()
            } else {
                $cell.Value2 = $val
            }
        } else {
            # every other column is stored as text in this workbook
            $cell.NumberFormat = "@"
            if ($null -eq $val) {
                $cell.ClearContents()
                $cell.NumberFormat = "@"
            } else {
                $cell.Value2 = $val
            }
        }
    }
}

function Apply-Permutation {
    param($ws, [int[]]$rows, [int[]]$sourceForTarget)
    # $rows[$i] gets the content that currently lives in row $sourceForTarget[$i]
    $buffers = @{}
    foreach ($r in $rows) {
        $buffers[$r] = Get-RowData $ws $r
    }
    for ($i = 0; $i -lt $rows.Length; $i++) {
        Set-RowData $ws $rows[$i] $buffers[$sourceForTarget[$i]]
    }
}

# Group 1: rows 93-94 swap
Apply-Permutation $ws @(93, 94) @(94, 93)

# Group 2: rows 165-166 swap
Apply-Permutation $ws @(165, 166) @(166, 165)

# Group 3: rows 262-264 rotate (262<-263, 263<-264, 264<-262)
Apply-Permutation $ws @(262, 263, 264) @(263, 264, 262)

# Group 4: rows 270, 271, 273 rotate (272 untouched / not part of the group)
Apply-Permutation $ws @(270, 271, 273) @(273, 270, 271)

# Group 5: rows 302-304 rotate (302<-304, 303<-302, 304<-303)
Apply-Permutation $ws @(302, 303, 304) @(304, 302, 303)

# Group 6: rows 345-348 rotate
Apply-Permutation $ws @(345, 346, 347, 348) @(346, 347, 348, 345)

# Every data row's timestamp gets updated to the new crawl time.
$tsRange = $ws.Range("O2:O409")
$tsRange.NumberFormat = "@"
$tsRange.Value2 = "2022-11-29 11:12:18"
